$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header row (13) --------------------------------------------------
# Mirrors the grey header style used by rows 3 / 8, but this command
# ("multiline") only needs firstname / lastname / birthday (no income),
# so only columns A-D get text while E13 just inherits the style.
$ws.Range("A3:E3").Copy()
$ws.Range("A13:E13").PasteSpecial(-4122)

$ws.Range("A13").Value = "COMMAND"
$ws.Range("B13").Value = "firstname"
$ws.Range("C13").Value = "lastname"
$ws.Range("D13").Value = "birthday"

# --- New data rows for "createPersonMultiline" ----------------------------
# Values are written in the order that first introduces each new shared
# string (createPersonMultiline, then Simon, Julian, Robert) so the
# resulting shared-string table matches the target ordering.
$ws.Range("A14").Value = "createPersonMultiline"
$ws.Range("C14").Value = "Kopp"

$ws.Range("B15").Value = "Simon"
$ws.Range("A15").Value = "createPersonMultiline"
$ws.Range("C15").Value = "Federer"

$ws.Range("B16").Value = "Julian"
$ws.Range("A16").Value = "createPersonMultiline"
$ws.Range("C16").Value = "Sallmer"

$ws.Range("B14").Value = "Robert"

# Birthdays, formatted as dates like the existing "birthday" columns.
$ws.Range("E4").Copy()
$ws.Range("D14:D16").PasteSpecial(-4122)

$ws.Range("D14").Value = 30266
$ws.Range("D15").Value = 41255
$ws.Range("D16").Value = 40544

# --- Cosmetics --------------------------------------------------------------
# Column A needs to be wide enough to fit "createPersonMultiline".
$ws.Columns.Item(1).ColumnWidth = 18.15

# Selection ends up on the (empty) E13 cell.
[void]$ws.Range("E13").Select()

# Page setup for printing (portrait, A4).
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
